# Correção das notas do fórum para matc65 em 2021.2
# For every student row (2..50) where column J ("nota_view") equals 4,
# reset all the daily-view columns (B..J) back to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 50 }

for ($r = 2; $r -le $lastRow; $r++) {
    $notaView = $ws.Cells.Item($r, 10).Value2  # column J
    if ($notaView -eq 4) {
        for ($c = 2; $c -le 10; $c++) {  # columns B..J
            $ws.Cells.Item($r, $c).Value = 0
        }
    }
}

$wb.Save()
